$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "answer"
$ws.Range("B1").Value = "question"
$ws.Range("D1").Value = "category"

$ws.Range("B1:D1").Font.Name = "Arial"

$ws.Range("D2").Select() | Out-Null
